$wb = $excel.ActiveWorkbook

# --- Rename sheets (new task-order timestamps) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16509960776488428"
$wb.Worksheets.Item(2).Name = "NB_TO-16509960794408443"
$wb.Worksheets.Item(3).Name = "RS_TO-16509960794408443"
$wb.Worksheets.Item(4).Name = "TOL_TO-1650996079488845"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1650996079568885"

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16509960776006942.csv"
$ws1.Range("B3").Value = "GNG_stims-16509960776246889.csv"
$ws1.Range("B4").Value = "go_stims-16509960776246889.csv"
$ws1.Range("B5").Value = "GNG_stims-16509960776488428.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16509960783368666.csv"
$ws2.Range("B3").Value = "ZB-match_4-16509960778728545.csv"
$ws2.Range("B4").Value = "ZB-match_9-16509960782168756.csv"
$ws2.Range("B5").Value = "TB-16509960794248435.csv"
$ws2.Range("B6").Value = "OB-16509960783128507.csv"
$ws2.Range("B7").Value = "TB-16509960788168538.csv"
$ws2.Range("B8").Value = "OB-16509960785448542.csv"
$ws2.Range("B9").Value = "TB-16509960788968437.csv"
$ws2.Range("B10").Value = "ZB-match_0-16509960776648495.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509960794568791.csv"
$ws4.Range("B3").Value = "ZM_stims-16509960794408443.csv"
$ws4.Range("B4").Value = "MM_stims-16509960794728796.csv"
$ws4.Range("B5").Value = "ZM_stims-16509960794568791.csv"
$ws4.Range("B6").Value = "MM_stims-1650996079488845.csv"
$ws4.Range("B7").Value = "ZM_stims-16509960794728796.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16509960795448773.csv"
$ws5.Range("B3").Value = "vSAT_stims-16509960795288796.csv"
$ws5.Range("B4").Value = "SAT_stims-16509960795128868.csv"
$ws5.Range("B5").Value = "SAT_stims-1650996079488845.csv"
